# case 1: update the A1:B3 values and narrow column A's width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A gets narrower (14.42578125 -> 13.7109375 "characters").
$ws.Columns.Item(1).ColumnWidth = 12.8

# Row 1
$ws.Range("A1").Value = 0.058184567309537812
$ws.Range("B1").Value = -0.058184567951516289

# Row 2
$ws.Range("A2").Value = 0.015642298518720947
$ws.Range("B2").Value = -0.015642299182742819

# Row 3
$ws.Range("A3").Value = -0.076691841547231257
$ws.Range("B3").Value = 0.076691840906616668

# Row 4 is left untouched.
